$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D header: "canonical SMILES"
$ws.Range("D2").Value = "canonical SMILES"

# Populate column D with the canonical SMILES values (duplicates of column C,
# since these microstates have no additional stereochemistry to encode)
$ws.Range("D3").Value = $ws.Range("C3").Value2
$ws.Range("D4").Value = $ws.Range("C4").Value2
$ws.Range("D5").Value = $ws.Range("C5").Value2
$ws.Range("D6").Value = $ws.Range("C6").Value2

# Set new column width (closest achievable value to the target 36.85546875
# given this runtime's pixel-quantized ColumnWidth conversion)
$ws.Columns.Item(4).ColumnWidth = 36
